$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRICES")

# Copy formatting (date style) from the row above into the new row's date cell
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New weekly price row
$ws.Range("A13").Value = 45633
$ws.Range("B13").Value = 22.2
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 2.5
$ws.Range("E13").Value = 11
$ws.Range("F13").Value = 4.2
$ws.Range("G13").Value = 5.7
$ws.Range("H13").Value = 1.8
$ws.Range("I13").Value = 1.64
$ws.Range("J13").Value = 54.95
$ws.Range("K13").Value = 23.95
$ws.Range("L13").Value = 419.95
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 104.95
$ws.Range("O13").Value = 67.5
$ws.Range("P13").Value = 42.12

# Restore the top-left view and move the selection as the author left it
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L17").Select()
